# Add method of adding lab member
# Replace the placeholder "aaa/bbb/ccc/ddd/eee" sample rows with a
# concrete example of three real lab members (sasaki, sato, tanaka),
# and drop the now-unused trailing rows (previously ddd/eee) so the
# sheet only keeps the header + 3 sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: aaa_position / a aa / aaa -> sasaki_position / 佐々木一郎 / sasaki
$ws.Range("A2").Value = "sasaki_position"
$ws.Range("B2").Value = "佐々木一郎"
$ws.Range("C2").Value = "sasaki"

# Row 3: bbb_position / b bb / bbb -> sato_position / 佐藤二郎 / sato
$ws.Range("A3").Value = "sato_position"
$ws.Range("B3").Value = "佐藤二郎"
$ws.Range("C3").Value = "sato"

# Row 4: ccc_position / c cc / ccc -> tanaka_position / 田中三郎 / tanaka
$ws.Range("A4").Value = "tanaka_position"
$ws.Range("B4").Value = "田中三郎"
$ws.Range("C4").Value = "tanaka"

# Remove the old rows 5 and 6 (ddd_position/d dd/ddd, eee_position/e ee/eee)
$ws.Rows("5:6").Delete()

# Move the selection, matching the saved view state in the edited file
[void]$ws.Range("A13").Select()
